$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) keeps Text formatting so numeric-looking
# strings (e.g. "1.00", "0.999") are stored exactly as text, matching
# the source data instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "90.541.26"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "3.135.67"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "215.81"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("D6").Value = "622.29"
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "1.14"
$ws.Range("E7").Value = "  +27.26%  "
$ws.Range("D8").Value = "0.363"
$ws.Range("E8").Value = "  -4.51%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").Value = "3.131.93"
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").Value = "0.737"
$ws.Range("E11").Value = "  +5.38%  "
$ws.Range("E12").Value = "  +5.87%  "
$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  -3.90%  "
$ws.Range("E14").Value = "  +4.85%  "
$ws.Range("D15").Value = "35.13"
$ws.Range("E15").Value = "  +6.13%  "
$ws.Range("D16").Value = "90.272.48"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("D17").Value = "3.711.91"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "3.138.87"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("D19").Value = "3.74"
$ws.Range("E19").Value = "  +3.88%  "
$ws.Range("D20").Value = "14.52"
$ws.Range("E20").Value = "  +5.04%  "
$ws.Range("D21").Value = "0.0000212"
$ws.Range("E21").Value = "  -7.73%  "
$ws.Range("D22").Value = "463.36"
$ws.Range("E22").Value = "  +7.45%  "
$ws.Range("D23").Value = "9.06"
$ws.Range("E23").Value = "  +5.92%  "
$ws.Range("E24").Value = "  +2.83%  "
$ws.Range("D25").Value = "5.92"
$ws.Range("E25").Value = "  +5.76%  "
$ws.Range("D26").Value = "94.97"
$ws.Range("E26").Value = "  +13.40%  "
$ws.Range("D27").Value = "12.26"
$ws.Range("E27").Value = "  +2.94%  "
$ws.Range("D28").Value = "3.309.31"
$ws.Range("E28").Value = "  +1.42%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -2.81%  "
$ws.Range("D31").Value = "9.19"
$ws.Range("E31").Value = "  +5.55%  "
$ws.Range("D32").Value = "0.214"
$ws.Range("E32").Value = "  +46.89%  "
$ws.Range("D33").Value = "26.72"
$ws.Range("E33").Value = "  +16.09%  "
$ws.Range("D34").Value = "517.50"
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("E36").Value = "  +5.04%  "
$ws.Range("D37").Value = "6.99"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("D39").Value = "3.59"
$ws.Range("E39").Value = "  -7.98%  "
$ws.Range("D40").Value = "0.0915"
$ws.Range("E40").Value = "  +27.89%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "0.426"
$ws.Range("E41").Value = "  +14.84%  "
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "22.21"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").Value = "0.752"
$ws.Range("E43").Value = "  -24.83%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "1.98"
$ws.Range("E45").Value = "  +5.31%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "0.727"
$ws.Range("E47").Value = "  +19.54%  "
$ws.Range("D48").Value = "4.73"
$ws.Range("E48").Value = "  +12.34%  "
$ws.Range("D49").Value = "150.34"
$ws.Range("E49").Value = "  +5.62%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "1.36"
$ws.Range("E50").Value = "  +7.86%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "45.32"
$ws.Range("E51").Value = "  +3.68%  "
